$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume figures (symbol list refresh).
# Cells are stored as text (e.g. "304.31", "3.84%"), so force
# text format before assigning to avoid Excel auto-converting
# the numeric-looking strings into real numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.84%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.61'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '13.50%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.091'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.26%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07824'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.99%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.251'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.31%'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.16%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '6.29%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9268'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.44%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09700'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.26%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1819'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.79%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08704'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03417'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.92%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09938'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.05%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001491'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.66%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005737'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.35%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.484'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.27%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.148'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.15%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3420'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.17%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.49%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.549'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '10.45%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04683'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.26%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001242'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.83%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004548'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.57%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001299'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.03%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002697'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-20.50%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01759'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '6.95%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04721'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.75%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007959'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '6.67%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1422'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.40%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008006'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-18.61%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002298'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.55%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009132'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.99%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006198'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.54%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.670'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '113.62%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002688'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '34.48%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001998'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
